$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.99658792072163393527262087
$ws.Range("E2").Value = 0.99658792072163393527262087

$ws.Range("D3").Value = 0.22421295730032750337379355
$ws.Range("E3").Value = 0.22421295730032750337379355

$ws.Range("D4").Value = 0.000011097559722188630489244851
$ws.Range("E4").Value = 0.000011097559722188630489244851

$ws.Range("D5").Value = 0.0000000000031616615676555430910369521
$ws.Range("E5").Value = 0.0000000000031616615676555430910369521

$ws.Range("D6").Value = 0.000000000000000000000000012231054588550840559670175
$ws.Range("E6").Value = 0.000000000000000000000000012231054588550840559670175

$ws.Range("D7").Value = 0.00059963981391389906995537107
$ws.Range("E7").Value = 0.99940036018608613410663111

$ws.Range("D9").Value = 0.99656463654046467670610809
$ws.Range("E9").Value = 0.0034353634595353228602110374

$ws.Range("D10").Value = 0.0012603299751487930142890281
$ws.Range("E10").Value = 0.99873967002485120048049794

$ws.Range("D11").Value = 0.000000000034363284491152822649335115
$ws.Range("E11").Value = 0.99999999996563671000870954
$ws.Range("F11").Value = 4.4127359390258789062500000
